$wb = $excel.ActiveWorkbook

# "Generate Report for Handback" - refresh the handoff/handback timestamps
# recorded for the primary (non content-duplicate) source file in each
# locale sheet, and roll the corresponding "Latest HO Xliff Generate Date"
# on the Overview sheet.

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-09-02 08:55:58"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-09-02 08:55:53"
$zhcn.Range("K2").Value = "2016-09-02 08:56:25"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-09-02 08:55:58"
$dede.Range("K2").Value = "2016-09-02 08:56:33"
